$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-05-27"

# Update the row label for May to reflect the new "through" date
$ws.Range("A6").Value = "May (through 05-27)"

# Update May row (row 6) values
$ws.Range("C6").Value = 40
$ws.Range("D6").Value = 53
$ws.Range("E6").Value = 44
$ws.Range("F6").Value = 40
$ws.Range("G6").Value = 56
$ws.Range("H6").Value = 101
$ws.Range("I6").Value = 96

# Update Total row (row 7) values
$ws.Range("C7").Value = 202
$ws.Range("D7").Value = 306
$ws.Range("E7").Value = 290
$ws.Range("F7").Value = 195
$ws.Range("G7").Value = 318
$ws.Range("H7").Value = 624
$ws.Range("I7").Value = 647
